$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Metadata sheet: bump the "Date" value to reflect the re-generation date.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# ---------------------------------------------------------------------------
# 2) Elements sheet: the two "Mapping" columns (AK = RIM Mapping,
#    AL = Spécification métier) are swapped so that the new business-mapping
#    column now comes first (AK) and the RIM mapping moves to AL.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Find the last used row on the Elements sheet so every row gets swapped.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $ws.Cells.Item($r, 37)  # column AK
    $alCell = $ws.Cells.Item($r, 38)  # column AL

    $akVal = $akCell.Value()
    $alVal = $alCell.Value()

    # Only touch cells whose value actually changes as a result of the swap,
    # so untouched rows (where both columns already hold the same content,
    # e.g. two blanks) are left byte-for-byte alone.
    if ($akVal -ne $alVal) {
        $akCell.Value = $alVal
        $alCell.Value = $akVal
    }
}

# Swap the column widths too: AK becomes the wide column (was AL's width),
# AL becomes the narrow column (was AK's width).
$ws.Columns.Item(37).ColumnWidth = 96.5
$ws.Columns.Item(38).ColumnWidth = 24.166666666666668
